# Auto-generated Excel COM-interop script applying the Goblin_Profits market-data refresh.
# For each changed row, update the currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N)
# to the newly scraped values. Some rows gain or lose a LeveProfit cell entirely (sparse data).
$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC!row 17
$ws_ALC.Range("H17").Value = 995.2889
$ws_ALC.Range("J17").Value = 1013.61365
$ws_ALC.Range("L17").Value = 3040.84095
$ws_ALC.Range("N17").Value = -3376.84095

# ALC!row 28
$ws_ALC.Range("H28").Value = 4959.609
$ws_ALC.Range("I28").Value = 5824
$ws_ALC.Range("K28").Value = 5824
$ws_ALC.Range("M28").Value = -5339

# ALC!row 33
$ws_ALC.Range("H33").Value = 1199.15
$ws_ALC.Range("I33").Value = 283.83334
$ws_ALC.Range("J33").Value = 2572.125
$ws_ALC.Range("K33").Value = 283.83334
$ws_ALC.Range("L33").Value = 2572.125
$ws_ALC.Range("M33").Value = -54.83334000000002
$ws_ALC.Range("N33").Value = -3030.125

# ALC!row 63
$ws_ALC.Range("H63").Value = 75000
$ws_ALC.Range("J63").Value = 75000
$ws_ALC.Range("L63").Value = 75000
$ws_ALC.Range("N63").Value = -76248

# ALC!row 66
$ws_ALC.Range("H66").Value = 75000
$ws_ALC.Range("J66").Value = 75000
$ws_ALC.Range("L66").Value = 225000
$ws_ALC.Range("N66").Value = -231240

# ALC!row 75
$ws_ALC.Range("H75").Value = 44999.5
$ws_ALC.Range("J75").Value = 44999.5
$ws_ALC.Range("L75").Value = 44999.5
$ws_ALC.Range("N75").Value = -46871.5

# ALC!row 78
$ws_ALC.Range("H78").Value = 44999.5
$ws_ALC.Range("J78").Value = 44999.5
$ws_ALC.Range("L78").Value = 134998.5
$ws_ALC.Range("N78").Value = -144358.5

# ALC!row 113
$ws_ALC.Range("H113").Value = 3281.5
$ws_ALC.Range("I113").Value = 3058
$ws_ALC.Range("K113").Value = 3058
$ws_ALC.Range("M113").Value = 196

# ALC!row 129
$ws_ALC.Range("H129").Value = 1613.9048
$ws_ALC.Range("J129").Value = 2124.7273
$ws_ALC.Range("L129").Value = 6374.1819
$ws_ALC.Range("N129").Value = -16374.1819

# ALC!row 137
$ws_ALC.Range("H137").Value = 14887.625
$ws_ALC.Range("I137").Value = 26651
$ws_ALC.Range("J137").Value = 3124.25
$ws_ALC.Range("K137").Value = 79953
$ws_ALC.Range("L137").Value = 9372.75
$ws_ALC.Range("M137").Value = -77403
$ws_ALC.Range("N137").Value = -14472.75

# ARM!row 32
$ws_ARM.Range("H32").Value = 18816.5
$ws_ARM.Range("I32").Value = 21379.8
$ws_ARM.Range("K32").Value = 21379.8
$ws_ARM.Range("M32").Value = -21092.8

# ARM!row 122
$ws_ARM.Range("H122").Value = 1791.6046
$ws_ARM.Range("I122").Value = 1528.8857
$ws_ARM.Range("K122").Value = 4586.6571
$ws_ARM.Range("M122").Value = -2136.6571

# BSM!row 86
$ws_BSM.Range("H86").Value = 1194.6111
$ws_BSM.Range("I86").Value = 1236.6471
$ws_BSM.Range("J86").Value = 480
$ws_BSM.Range("K86").Value = 1236.6471
$ws_BSM.Range("L86").Value = 480
$ws_BSM.Range("M86").Value = -113.6470999999999
$ws_BSM.Range("N86").Value = -2726

# BSM!row 89
$ws_BSM.Range("H89").Value = 1194.6111
$ws_BSM.Range("I89").Value = 1236.6471
$ws_BSM.Range("J89").Value = 480
$ws_BSM.Range("K89").Value = 6183.2355
$ws_BSM.Range("L89").Value = 2400
$ws_BSM.Range("M89").Value = -567.2354999999998
$ws_BSM.Range("N89").Value = -13632

# BSM!row 92
$ws_BSM.Range("H92").Value = 49133.332
$ws_BSM.Range("I92").Value = 48000
$ws_BSM.Range("J92").Value = 49700
$ws_BSM.Range("K92").Value = 48000
$ws_BSM.Range("L92").Value = 49700
$ws_BSM.Range("M92").Value = -45504
$ws_BSM.Range("N92").Value = -54692

# BSM!row 107
$ws_BSM.Range("H107").Value = 3358.068
$ws_BSM.Range("I107").Value = 1284.862
$ws_BSM.Range("K107").Value = 1284.862
$ws_BSM.Range("M107").Value = 635.1379999999999

# BSM!row 134
$ws_BSM.Range("H134").Value = 2693.2104
$ws_BSM.Range("I134").Value = 2628.5386
$ws_BSM.Range("K134").Value = 7885.6158
$ws_BSM.Range("M134").Value = -5350.6158

# CRP!row 22
$ws_CRP.Range("H22").Value = 1130.909
$ws_CRP.Range("I22").Value = 816.9048
$ws_CRP.Range("J22").Value = 1680.4166
$ws_CRP.Range("K22").Value = 816.9048
$ws_CRP.Range("L22").Value = 1680.4166
$ws_CRP.Range("M22").Value = -466.9048
$ws_CRP.Range("N22").Value = -2380.4166

# CRP!row 58
$ws_CRP.Range("H58").Value = 1533.421
$ws_CRP.Range("I58").Value = 1608.7273
$ws_CRP.Range("J58").Value = 1429.875
$ws_CRP.Range("K58").Value = 1608.7273
$ws_CRP.Range("L58").Value = 1429.875
$ws_CRP.Range("M58").Value = -1405.7273
$ws_CRP.Range("N58").Value = -1835.875

# CRP!row 132
$ws_CRP.Range("H132").Value = 1856.2858
$ws_CRP.Range("I132").Value = 1822.2354
$ws_CRP.Range("K132").Value = 5466.706200000001
$ws_CRP.Range("M132").Value = -2936.706200000001

# CRP!row 136
$ws_CRP.Range("H136").Value = 1533.421
$ws_CRP.Range("I136").Value = 1608.7273
$ws_CRP.Range("J136").Value = 1429.875
$ws_CRP.Range("K136").Value = 4826.1819
$ws_CRP.Range("L136").Value = 4289.625
$ws_CRP.Range("M136").Value = -2276.1819
$ws_CRP.Range("N136").Value = -9389.625

# CUL!row 92
$ws_CUL.Range("H92").Value = 2288.3333
$ws_CUL.Range("J92").Value = 1340.6
$ws_CUL.Range("L92").Value = 4021.8
$ws_CUL.Range("N92").Value = -6517.799999999999

# CUL!row 113
$ws_CUL.Range("H113").Value = 1503.9166
$ws_CUL.Range("J113").Value = 1766.3334
$ws_CUL.Range("L113").Value = 5299.0002
$ws_CUL.Range("N113").Value = -9639.0002

# CUL!row 140
$ws_CUL.Range("H140").Value = 64206.188
$ws_CUL.Range("I140").Value = 78028.92
$ws_CUL.Range("K140").Value = 234086.76
$ws_CUL.Range("M140").Value = -228906.76

# GSM!row 107
$ws_GSM.Range("H107").Value = 555.8
$ws_GSM.Range("I107").Value = 484.22223
$ws_GSM.Range("K107").Value = 484.22223
$ws_GSM.Range("M107").Value = 1435.77777

# GSM!row 132
$ws_GSM.Range("H132").Value = 2580.5
$ws_GSM.Range("I132").Value = 2672.2632
$ws_GSM.Range("K132").Value = 8016.7896
$ws_GSM.Range("M132").Value = -5486.7896

# LTW!row 36
$ws_LTW.Range("H36").Value = 50000
$ws_LTW.Range("J36").Value = 50000
$ws_LTW.Range("L36").Value = 50000
$ws_LTW.Range("N36").Value = -51124

# LTW!row 132
$ws_LTW.Range("H132").Value = 5188.4707
$ws_LTW.Range("I132").Value = 5099.875
$ws_LTW.Range("J132").Value = 5267.222
$ws_LTW.Range("K132").Value = 15299.625
$ws_LTW.Range("L132").Value = 15801.666
$ws_LTW.Range("M132").Value = -12769.625
$ws_LTW.Range("N132").Value = -20861.666

# LTW!row 136
$ws_LTW.Range("H136").Value = 5399.9
$ws_LTW.Range("I136").Value = 5333.222
$ws_LTW.Range("K136").Value = 15999.666
$ws_LTW.Range("M136").Value = -13449.666

# WVR!row 62
$ws_WVR.Range("H62").Value = 10957.143
$ws_WVR.Range("I62").Value = 5366.6665
$ws_WVR.Range("J62").Value = 12481.818
$ws_WVR.Range("K62").Value = 5366.6665
$ws_WVR.Range("L62").Value = 12481.818
$ws_WVR.Range("M62").Value = -4742.6665
$ws_WVR.Range("N62").Value = -13729.818

# WVR!row 65
$ws_WVR.Range("H65").Value = 10957.143
$ws_WVR.Range("I65").Value = 5366.6665
$ws_WVR.Range("J65").Value = 12481.818
$ws_WVR.Range("K65").Value = 26833.3325
$ws_WVR.Range("L65").Value = 62409.09
$ws_WVR.Range("M65").Value = -23713.3325
$ws_WVR.Range("N65").Value = -68649.09

# WVR!row 81
$ws_WVR.Range("H81").Value = 2599.8
$ws_WVR.Range("I81").Value = 0
$ws_WVR.Range("J81").Value = 2599.8
$ws_WVR.Range("K81").Value = 0
$ws_WVR.Range("L81").Value = 5199.6
$ws_WVR.Range("M81").ClearContents()
$ws_WVR.Range("N81").Value = -7321.6

# WVR!row 84
$ws_WVR.Range("H84").Value = 2599.8
$ws_WVR.Range("I84").Value = 0
$ws_WVR.Range("J84").Value = 2599.8
$ws_WVR.Range("K84").Value = 0
$ws_WVR.Range("L84").Value = 25998
$ws_WVR.Range("M84").ClearContents()
$ws_WVR.Range("N84").Value = -36606

# WVR!row 107
$ws_WVR.Range("H107").Value = 3065.9033
$ws_WVR.Range("I107").Value = 1360.1923
$ws_WVR.Range("K107").Value = 4080.5769
$ws_WVR.Range("M107").Value = -2160.5769
